$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.976.67"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "3.247.45"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "396.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.583"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.44%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0962"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.75%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "3.755.91"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "3.233.57"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  -3.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.97%  "
$ws.Range("D19").Value = "56.831.77"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("E21").Value = "  +9.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "293.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.33%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "41.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "136.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("E47").Value = "  +8.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").Value = "2.153.88"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  -5.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.80%  "
